$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 30 de Marzo de 2020 a las 17:20"

# Row 4
$ws.Cells.Item(4, 2).Value = 145088
$ws.Cells.Item(4, 3).Value = 1597
$ws.Cells.Item(4, 4).Value = 4573
$ws.Cells.Item(4, 5).Value = 137910
$ws.Cells.Item(4, 6).Value = 2970
$ws.Cells.Item(4, 7).Value = 23
$ws.Cells.Item(4, 8).Value = 2605

# Row 18
$ws.Cells.Item(18, 1).Value = "Canada"
$ws.Cells.Item(18, 2).Value = 6671
$ws.Cells.Item(18, 3).Value = 351
$ws.Cells.Item(18, 4).Value = 1014
$ws.Cells.Item(18, 5).Value = 5590
$ws.Cells.Item(18, 6).Value = 120
$ws.Cells.Item(18, 7).Value = 2
$ws.Cells.Item(18, 8).Value = 67

# Row 19
$ws.Cells.Item(19, 1).Value = "Portugal"
$ws.Cells.Item(19, 2).Value = 6408
$ws.Cells.Item(19, 3).Value = 446
$ws.Cells.Item(19, 4).Value = 43
$ws.Cells.Item(19, 5).Value = 6225
$ws.Cells.Item(19, 6).Value = 164
$ws.Cells.Item(19, 7).Value = 21
$ws.Cells.Item(19, 8).Value = 140

# Row 25
$ws.Cells.Item(25, 2).Value = 2896
$ws.Cells.Item(25, 3).Value = 79
$ws.Cells.Item(25, 4).Value = 11
$ws.Cells.Item(25, 5).Value = 2868
$ws.Cells.Item(25, 6).Value = 52
$ws.Cells.Item(25, 7).Value = 1
$ws.Cells.Item(25, 8).Value = 17

# Row 43
$ws.Cells.Item(43, 2).Value = 1212
$ws.Cells.Item(43, 3).Value = 56
$ws.Cells.Item(43, 4).Value = 52
$ws.Cells.Item(43, 5).Value = 1117
$ws.Cells.Item(43, 6).Value = 66
$ws.Cells.Item(43, 7).Value = 4
$ws.Cells.Item(43, 8).Value = 43

# Row 71
$ws.Cells.Item(71, 2).Value = 446
$ws.Cells.Item(71, 3).Value = 8
$ws.Cells.Item(71, 4).Value = 35
$ws.Cells.Item(71, 5).Value = 400
$ws.Cells.Item(71, 6).Value = 3
$ws.Cells.Item(71, 7).Value = 1
$ws.Cells.Item(71, 8).Value = 11

# Row 73
$ws.Cells.Item(73, 2).Value = 359
$ws.Cells.Item(73, 3).Value = 13
$ws.Cells.Item(73, 4).Value = 17
$ws.Cells.Item(73, 5).Value = 334
$ws.Cells.Item(73, 6).Value = 13
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 8

# Row 85
$ws.Cells.Item(85, 2).Value = 263
$ws.Cells.Item(85, 3).Value = 0
$ws.Cells.Item(85, 4).Value = 13
$ws.Cells.Item(85, 5).Value = 248
$ws.Cells.Item(85, 6).Value = 44
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 2

# Row 87
$ws.Cells.Item(87, 1).Value = "Republica de Chipre"
$ws.Cells.Item(87, 2).Value = 230
$ws.Cells.Item(87, 3).Value = 16
$ws.Cells.Item(87, 4).Value = 15
$ws.Cells.Item(87, 5).Value = 208
$ws.Cells.Item(87, 6).Value = 3
$ws.Cells.Item(87, 7).Value = 1
$ws.Cells.Item(87, 8).Value = 7

# Row 88
$ws.Cells.Item(88, 1).Value = "San Marino"
$ws.Cells.Item(88, 2).Value = 230
$ws.Cells.Item(88, 3).Value = 6
$ws.Cells.Item(88, 4).Value = 13
$ws.Cells.Item(88, 5).Value = 192
$ws.Cells.Item(88, 6).Value = 16
$ws.Cells.Item(88, 7).Value = 3
$ws.Cells.Item(88, 8).Value = 25

# Row 89
$ws.Cells.Item(89, 1).Value = "Albania"
$ws.Cells.Item(89, 2).Value = 223
$ws.Cells.Item(89, 3).Value = 11
$ws.Cells.Item(89, 4).Value = 44
$ws.Cells.Item(89, 5).Value = 168
$ws.Cells.Item(89, 6).Value = 7
$ws.Cells.Item(89, 7).Value = 1
$ws.Cells.Item(89, 8).Value = 11

# Row 90
$ws.Cells.Item(90, 1).Value = "Burkina Faso"
$ws.Cells.Item(90, 2).Value = 222
$ws.Cells.Item(90, 3).Value = 0
$ws.Cells.Item(90, 4).Value = 23
$ws.Cells.Item(90, 5).Value = 187
$ws.Cells.Item(90, 6).Value = 0
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 12

# Row 101
$ws.Cells.Item(101, 1).Value = "Afganistan"
$ws.Cells.Item(101, 2).Value = 145
$ws.Cells.Item(101, 3).Value = 25
$ws.Cells.Item(101, 4).Value = 2
$ws.Cells.Item(101, 5).Value = 139
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 4

# Row 102
$ws.Cells.Item(102, 1).Value = "Honduras"
$ws.Cells.Item(102, 2).Value = 139
$ws.Cells.Item(102, 3).Value = 29
$ws.Cells.Item(102, 4).Value = 3
$ws.Cells.Item(102, 5).Value = 133
$ws.Cells.Item(102, 6).Value = 4
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 3

# Row 103
$ws.Cells.Item(103, 1).Value = "Cuba"
$ws.Cells.Item(103, 2).Value = 139
$ws.Cells.Item(103, 3).Value = 0
$ws.Cells.Item(103, 4).Value = 4
$ws.Cells.Item(103, 5).Value = 132
$ws.Cells.Item(103, 6).Value = 2
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 3

# Row 104
$ws.Cells.Item(104, 1).Value = "Camerun"
$ws.Cells.Item(104, 2).Value = 139
$ws.Cells.Item(104, 3).Value = 0
$ws.Cells.Item(104, 4).Value = 5
$ws.Cells.Item(104, 5).Value = 128
$ws.Cells.Item(104, 6).Value = 0
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 6

# Row 105
$ws.Cells.Item(105, 1).Value = "Venezuela"
$ws.Cells.Item(105, 2).Value = 129
$ws.Cells.Item(105, 3).Value = 10
$ws.Cells.Item(105, 4).Value = 39
$ws.Cells.Item(105, 5).Value = 87
$ws.Cells.Item(105, 6).Value = 6
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = 3

# Row 106
$ws.Cells.Item(106, 1).Value = "Mauricio"
$ws.Cells.Item(106, 2).Value = 128
$ws.Cells.Item(106, 3).Value = 21
$ws.Cells.Item(106, 4).Value = 0
$ws.Cells.Item(106, 5).Value = 125
$ws.Cells.Item(106, 6).Value = 1
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 3

# Row 107
$ws.Cells.Item(107, 1).Value = "Brunei"
$ws.Cells.Item(107, 2).Value = 127
$ws.Cells.Item(107, 3).Value = 1
$ws.Cells.Item(107, 4).Value = 38
$ws.Cells.Item(107, 5).Value = 88
$ws.Cells.Item(107, 6).Value = 3
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 1

# Row 108
$ws.Cells.Item(108, 1).Value = "Sri Lanka"
$ws.Cells.Item(108, 2).Value = 122
$ws.Cells.Item(108, 3).Value = 5
$ws.Cells.Item(108, 4).Value = 15
$ws.Cells.Item(108, 5).Value = 105
$ws.Cells.Item(108, 6).Value = 5
$ws.Cells.Item(108, 7).Value = 1
$ws.Cells.Item(108, 8).Value = 2

# Row 109
$ws.Cells.Item(109, 1).Value = "Estado de Palestina"
$ws.Cells.Item(109, 2).Value = 115
$ws.Cells.Item(109, 3).Value = 6
$ws.Cells.Item(109, 4).Value = 18
$ws.Cells.Item(109, 5).Value = 96
$ws.Cells.Item(109, 6).Value = 0
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 1

# Row 110
$ws.Cells.Item(110, 1).Value = "Nigeria"
$ws.Cells.Item(110, 2).Value = 111
$ws.Cells.Item(110, 3).Value = 0
$ws.Cells.Item(110, 4).Value = 3
$ws.Cells.Item(110, 5).Value = 107
$ws.Cells.Item(110, 6).Value = 0
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 1

# Row 114
$ws.Cells.Item(114, 2).Value = 97
$ws.Cells.Item(114, 3).Value = 16
$ws.Cells.Item(114, 4).Value = 0
$ws.Cells.Item(114, 5).Value = 93
$ws.Cells.Item(114, 6).Value = 3
$ws.Cells.Item(114, 7).Value = 3
$ws.Cells.Item(114, 8).Value = 4

# Row 127
$ws.Cells.Item(127, 1).Value = "Isla de Man"
$ws.Cells.Item(127, 2).Value = 49
$ws.Cells.Item(127, 3).Value = 7
$ws.Cells.Item(127, 4).Value = 0
$ws.Cells.Item(127, 5).Value = 49
$ws.Cells.Item(127, 6).Value = 0
$ws.Cells.Item(127, 7).Value = 0
$ws.Cells.Item(127, 8).Value = 0

# Row 128
$ws.Cells.Item(128, 1).Value = "Banglades"
$ws.Cells.Item(128, 2).Value = 49
$ws.Cells.Item(128, 3).Value = 1
$ws.Cells.Item(128, 4).Value = 19
$ws.Cells.Item(128, 5).Value = 25
$ws.Cells.Item(128, 6).Value = 1
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 5

# Row 134
$ws.Cells.Item(134, 1).Value = "Jamaica"
$ws.Cells.Item(134, 2).Value = 36
$ws.Cells.Item(134, 3).Value = 4
$ws.Cells.Item(134, 4).Value = 2
$ws.Cells.Item(134, 5).Value = 33
$ws.Cells.Item(134, 6).Value = 0
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 1

# Row 135
$ws.Cells.Item(135, 1).Value = "Guatemala"
$ws.Cells.Item(135, 2).Value = 36
$ws.Cells.Item(135, 3).Value = 2
$ws.Cells.Item(135, 4).Value = 10
$ws.Cells.Item(135, 5).Value = 25
$ws.Cells.Item(135, 6).Value = 1
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 1

# Row 137
$ws.Cells.Item(137, 1).Value = "Zambia"
$ws.Cells.Item(137, 2).Value = 35
$ws.Cells.Item(137, 3).Value = 6
$ws.Cells.Item(137, 4).Value = 0
$ws.Cells.Item(137, 5).Value = 35
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 0

# Row 177
$ws.Cells.Item(177, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(177, 2).Value = 7
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 0
$ws.Cells.Item(177, 5).Value = 7
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 0

# Row 178
$ws.Cells.Item(178, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(178, 2).Value = 7
$ws.Cells.Item(178, 3).Value = 5
$ws.Cells.Item(178, 4).Value = 0
$ws.Cells.Item(178, 5).Value = 7
$ws.Cells.Item(178, 6).Value = 0
$ws.Cells.Item(178, 7).Value = 0
$ws.Cells.Item(178, 8).Value = 0

# Row 179
$ws.Cells.Item(179, 1).Value = "Zimbabue"

# Row 180
$ws.Cells.Item(180, 1).Value = "Gabon"

# Row 182
$ws.Cells.Item(182, 1).Value = "Benin"

# Row 183
$ws.Cells.Item(183, 1).Value = "Santa Sede"

# Row 184
$ws.Cells.Item(184, 1).Value = "San Martin (Parte Holandesa)"

# Row 187
$ws.Cells.Item(187, 1).Value = "Montserrat"

# Row 188
$ws.Cells.Item(188, 1).Value = "San Bartolome"

# Row 189
$ws.Cells.Item(189, 1).Value = "Republica del Chad"
$ws.Cells.Item(189, 2).Value = 5
$ws.Cells.Item(189, 3).Value = 2
$ws.Cells.Item(189, 4).Value = 0
$ws.Cells.Item(189, 5).Value = 5
$ws.Cells.Item(189, 6).Value = 0
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 0

# Row 190
$ws.Cells.Item(190, 1).Value = "Fiyi"
$ws.Cells.Item(190, 2).Value = 5
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(190, 4).Value = 0
$ws.Cells.Item(190, 5).Value = 5
$ws.Cells.Item(190, 6).Value = 0
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 0

# Row 191
$ws.Cells.Item(191, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(191, 2).Value = 5
$ws.Cells.Item(191, 3).Value = 1
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 5).Value = 5
$ws.Cells.Item(191, 6).Value = 0
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 0

# Row 192
$ws.Cells.Item(192, 1).Value = "Nepal"
$ws.Cells.Item(192, 2).Value = 5
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 1
$ws.Cells.Item(192, 5).Value = 4
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 0

# Row 193
$ws.Cells.Item(193, 1).Value = "Mauritania"
$ws.Cells.Item(193, 2).Value = 5
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 2
$ws.Cells.Item(193, 5).Value = 3
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 0

# Row 194
$ws.Cells.Item(194, 1).Value = "Butan"
